$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 61, shifting existing rows 61:87 down to 62:88.
$ws.Rows(61).Insert()

# Populate the newly inserted row 61 with the weekly price-report entry.
$ws.Range("A61").Value = 4
$ws.Range("B61").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C61").Value = "Los Lagos"
$ws.Range("D61").Value = 44523
$ws.Range("E61").Value = 10
$ws.Range("F61").Value = 100112022
$ws.Range("G61").Value = "Arveja Verde"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 150
$ws.Range("K61").Value = 18000
$ws.Range("L61").Value = 18000
$ws.Range("M61").Value = 18000
$ws.Range("N61").Value = "$/saco 25 kilos"
$ws.Range("O61").Value = "Región del Maule"
$ws.Range("P61").Value = 720
$ws.Range("Q61").Value = 25
$ws.Range("R61").Value = "Hortaliza"
